# Add new test case row (TEST_18) to "Hoja 1" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja 1")

$row = 18

$ws.Cells.Item($row, 1).Value  = "ENTRADA"
$ws.Cells.Item($row, 2).Value  = "EC"
$ws.Cells.Item($row, 3).Value  = "fichero"
$ws.Cells.Item($row, 4).Value  = "INVALID"
$ws.Cells.Item($row, 5).Value  = "TEST_18"
$ws.Cells.Item($row, 6).Value  = "Fichero facío"
$ws.Cells.Item($row, 7).Value  = "Comprobar si hay datos en el fichero"
$ws.Cells.Item($row, 8).Value  = "NONE"
$ws.Cells.Item($row, 9).Value  = "NONE"
$ws.Cells.Item($row, 10).Value = "NONE"
$ws.Cells.Item($row, 11).Value = "ERROR"

# Copy the plain data-cell formatting (style used by A17, shared by every
# column on row 17 except J) across the whole new row, including J18 -
# unlike row 17, J18 holds a normal "NONE" value, not a long REASON string.
$ws.Range("A17").Copy()
$ws.Range("A18:K18").PasteSpecial(-4122)  # xlPasteFormats

# Match row 17/18's taller, manually-set row height.
$ws.Rows.Item(18).RowHeight = 15.75

$ws.Range("A18").Select()
